$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 24; $row++) {
    $ws.Range("L$row").Formula = "=K$row/3600"
    $ws.Range("M$row").Formula = "=J$row-L$row"
    $ws.Range("P$row").Formula = "=J$row/L$row%"
}
